$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (font/border/alignment) from the existing H1 header cell
# onto the two new header cells I1 and J1 before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for column I (I0)
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 3

# New data values for column J (IF)
$ws.Range("J2").Value = 5
$ws.Range("J3").Value = 9
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 4
